$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 3 (the old January row with blank species and J=0),
# which shifts rows 4-7 up to become rows 3-6.
$ws.Rows.Item(3).Delete()
